$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2-11 changes from serial date 45188 (2023-09-19)
# to serial date 45189 (2023-09-20).
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45189
}
